$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain a text/string value even when the
    # content looks like a number (matches source data stored as text).
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "42.614.74"
$ws.Range("E2").Value = "  -6.85%  "

# Row 3
$ws.Range("D3").Value = "2.216.73"
$ws.Range("E3").Value = "  -7.55%  "

# Row 4
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
Set-TextValue $ws.Range("D5") "312.58"
$ws.Range("E5").Value = "  -1.91%  "

# Row 6
Set-TextValue $ws.Range("D6") "97.49"
$ws.Range("E6").Value = "  -14.26%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.570"
$ws.Range("E7").Value = "  -10.15%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.555"
$ws.Range("E9").Value = "  -11.36%  "

# Row 10
Set-TextValue $ws.Range("D10") "36.70"
$ws.Range("E10").Value = "  -12.43%  "

# Row 11
Set-TextValue $ws.Range("D11") "53.46"
$ws.Range("E11").Value = "  -4.79%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0829"
$ws.Range("E12").Value = "  -10.78%  "

# Row 13
Set-TextValue $ws.Range("D13") "7.57"
$ws.Range("E13").Value = "  -13.17%  "

# Row 14
$ws.Range("E14").Value = "  -4.82%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.867"
$ws.Range("E15").Value = "  -13.79%  "

# Row 16
$ws.Range("D16").Value = "2.554.05"
$ws.Range("E16").Value = "  -7.56%  "

# Row 17
Set-TextValue $ws.Range("D17") "13.85"
$ws.Range("E17").Value = "  -12.46%  "

# Row 18
$ws.Range("D18").Value = "2.211.80"
$ws.Range("E18").Value = "  -7.70%  "

# Row 19
$ws.Range("D19").Value = "42.543.33"
$ws.Range("E19").Value = "  -6.89%  "

# Row 20
Set-TextValue $ws.Range("D20") "13.78"
$ws.Range("E20").Value = "  +2.38%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.60"
$ws.Range("E21").Value = "  -11.96%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0944"
$ws.Range("E22").Value = "  -12.93%  "

# Row 23
Set-TextValue $ws.Range("D23") "3.25"
$ws.Range("E23").Value = "  -8.27%  "

# Row 24
Set-TextValue $ws.Range("D24") "64.46"
$ws.Range("E24").Value = "  -13.74%  "

# Row 25
Set-TextValue $ws.Range("D25") "233.38"
$ws.Range("E25").Value = "  -11.72%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.10"
$ws.Range("E26").Value = "  -10.64%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.995"
$ws.Range("E27").Value = "  -0.56%  "

# Row 28
Set-TextValue $ws.Range("D28") "10.10"
$ws.Range("E28").Value = "  -10.76%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D29") "2.17"
$ws.Range("E29").Value = "  -7.86%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D30") "6.48"
$ws.Range("E30").Value = "  -15.14%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.0881"
$ws.Range("E31").Value = "  -9.13%  "

# Row 32
Set-TextValue $ws.Range("D32") "20.31"
$ws.Range("E32").Value = "  -10.81%  "

# Row 33
Set-TextValue $ws.Range("D33") "157.35"
$ws.Range("E33").Value = "  -8.96%  "

# Row 34
Set-TextValue $ws.Range("D34") "32.47"
$ws.Range("E34").Value = "  -17.05%  "

# Row 35
Set-TextValue $ws.Range("D35") "2.71"
$ws.Range("E35").Value = "  -8.30%  "

# Row 36
$ws.Range("E36").Value = "  +0.92%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.122"
$ws.Range("E37").Value = "  -8.29%  "

# Row 38
Set-TextValue $ws.Range("D38") "4.40"
$ws.Range("E38").Value = "  -10.58%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.81"
$ws.Range("E39").Value = "  +2.05%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.104"
$ws.Range("E40").Value = "  -11.73%  "

# Row 41
Set-TextValue $ws.Range("D41") "3.51"
$ws.Range("E41").Value = "  -15.13%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.0320"
$ws.Range("E42").Value = "  -11.80%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.01"
$ws.Range("E43").Value = "  +0.41%  "

# Row 44
$ws.Range("D44").Value = "1.793.09"
$ws.Range("E44").Value = "  +7.66%  "

# Row 45
Set-TextValue $ws.Range("D45") "88.34"
$ws.Range("E45").Value = "  -12.54%  "

# Row 46
Set-TextValue $ws.Range("D46") "11.90"
$ws.Range("E46").Value = "  -11.86%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.205"
$ws.Range("E47").Value = "  -14.87%  "

# Row 48
Set-TextValue $ws.Range("D48") "76.95"
$ws.Range("E48").Value = "  -12.59%  "

# Row 49
Set-TextValue $ws.Range("D49") "5.35"
$ws.Range("E49").Value = "  -6.22%  "

# Row 50
Set-TextValue $ws.Range("D50") "59.92"
$ws.Range("E50").Value = "  -16.72%  "

# Row 51
Set-TextValue $ws.Range("D51") "8.54"
$ws.Range("E51").Value = "  -9.74%  "
